$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -9.361830160744564
$ws.Range("C2").Value = 79.9214580198452
$ws.Range("D2").Value = 8.109283971907509
$ws.Range("E2").Value = 0.00005449597197199146
$ws.Range("F2").Value = -1.336606832821891
$ws.Range("G2").Value = -0.543091955582276
$ws.Range("H2").Value = 0.2818183968919779
$ws.Range("I2").Value = 1.812345802972573
$ws.Range("J2").Value = 4.353799460192466
$ws.Range("K2").Value = 249
$ws.Range("L2").Value = -29.8862848846466
$ws.Range("N2").Value = 4.355584533945237
$ws.Range("O2").Value = 5.022251200611904

$ws.Range("B3").Value = 6.596534342844713
$ws.Range("C3").Value = 1.825852598549732
$ws.Range("D3").Value = 0.0001378031027128604
$ws.Range("E3").Value = 75.20885893170779
$ws.Range("F3").Value = 0.2961473385211364
$ws.Range("G3").Value = -0.6011903421046492
$ws.Range("H3").Value = 1.707038627490326
$ws.Range("I3").Value = -0.5890721322461243
$ws.Range("J3").Value = 4.353762573808641
$ws.Range("K3").Value = 353
$ws.Range("L3").Value = -24.69881904272587
$ws.Range("N3").Value = 4.35572411496149
$ws.Range("O3").Value = 5.022390781628157

$ws.Range("B4").Value = 0.0003578390485837374
$ws.Range("C4").Value = 26.48330213504671
$ws.Range("D4").Value = 74.18715757271076
$ws.Range("E4").Value = 2.768469704767545
$ws.Range("F4").Value = 1.455781544625709
$ws.Range("G4").Value = -0.009869048897695798
$ws.Range("H4").Value = -0.6828996761626946
$ws.Range("I4").Value = 0.4053679610471881
$ws.Range("J4").Value = 4.353781715690026
$ws.Range("K4").Value = 18
$ws.Range("L4").Value = -40.13340457506403
$ws.Range("N4").Value = 4.355733012259954
$ws.Range("O4").Value = 5.022399678926621

$ws.Range("B5").Value = 7.235222724324111
$ws.Range("C5").Value = 1.823592883824094
$ws.Range("D5").Value = 0.1732319419628537
$ws.Range("E5").Value = 74.46176761665646
$ws.Range("F5").Value = 0.2317875963304656
$ws.Range("G5").Value = -0.6868283319011741
$ws.Range("H5").Value = 0.6997131497254809
$ws.Range("I5").Value = -0.6233165677115948
$ws.Range("J5").Value = 4.353796303221486
$ws.Range("K5").Value = 136
$ws.Range("L5").Value = -23.11438027717823
$ws.Range("N5").Value = 4.3557415291715
$ws.Range("O5").Value = 5.022408195838167

$ws.Range("B6").Value = 2.375872132867816
$ws.Range("C6").Value = 0.0001071321012569115
$ws.Range("D6").Value = 55.21726756618601
$ws.Range("E6").Value = 71.76632244377979
$ws.Range("F6").Value = 0.4479931904071734
$ws.Range("G6").Value = 1.257083984185609
$ws.Range("H6").Value = -0.03676376659192426
$ws.Range("I6").Value = -0.7983122541075389
$ws.Range("J6").Value = 4.353808698294706
$ws.Range("K6").Value = 774
$ws.Range("L6").Value = -60.6001529040814
$ws.Range("N6").Value = 4.355747789085157
$ws.Range("O6").Value = 5.022414455751824

$ws.Range("B7").Value = 16.03804488700605
$ws.Range("C7").Value = 14.49656565995068
$ws.Range("D7").Value = 74.69710326583348
$ws.Range("E7").Value = 0.07144430964662618
$ws.Range("F7").Value = -1.547401700763939
$ws.Range("G7").Value = 0.1835092744198308
$ws.Range("H7").Value = -0.5382717719789614
$ws.Range("I7").Value = 0.8114934974099586
$ws.Range("J7").Value = 4.353806728819031
$ws.Range("K7").Value = 974
$ws.Range("L7").Value = -36.45841838058163
$ws.Range("N7").Value = 4.355757354098976
$ws.Range("O7").Value = 5.022424020765643

$ws.Range("B8").Value = 7.676241468781314
$ws.Range("C8").Value = 21.60751300598568
$ws.Range("D8").Value = 58.38320844267371
$ws.Range("E8").Value = 0.001275615317276173
$ws.Range("F8").Value = 0.2714940977360194
$ws.Range("G8").Value = -0.7794024848719248
$ws.Range("H8").Value = -0.5519685565960033
$ws.Range("I8").Value = 1.371028643870078
$ws.Range("J8").Value = 4.353758200196836
$ws.Range("K8").Value = 108
$ws.Range("L8").Value = -26.13370530623557
$ws.Range("N8").Value = 4.355764287152763
$ws.Range("O8").Value = 5.02243095381943

$ws.Range("B9").Value = 0.747428462616138
$ws.Range("C9").Value = 74.05165941048614
$ws.Range("D9").Value = 8.054987754254459
$ws.Range("E9").Value = 1.881876605834397
$ws.Range("F9").Value = 0.5379062400866643
$ws.Range("G9").Value = -0.7302051705165971
$ws.Range("H9").Value = -0.1036781222442318
$ws.Range("I9").Value = 0.3298165943394573
$ws.Range("J9").Value = 4.353799575730104
$ws.Range("K9").Value = 961
$ws.Range("L9").Value = -18.48984247227505
$ws.Range("N9").Value = 4.35581188612813
$ws.Range("O9").Value = 5.022478552794797

$ws.Range("B10").Value = 2.69765586013661
$ws.Range("C10").Value = 80.71012257875209
$ws.Range("D10").Value = 0.8766402467110316
$ws.Range("E10").Value = 11.89895577184867
$ws.Range("F10").Value = 0.1841065939362734
$ws.Range("G10").Value = -0.8541930096725487
$ws.Range("H10").Value = 0.5558024644633233
$ws.Range("I10").Value = -0.1453879663425757
$ws.Range("J10").Value = 4.353778366623544
$ws.Range("K10").Value = 726
$ws.Range("L10").Value = -18.74199211138344
$ws.Range("N10").Value = 4.355977461948342
$ws.Range("O10").Value = 5.022644128615009

$ws.Range("B11").Value = 0.0306258427404588
$ws.Range("C11").Value = 94.34855615777451
$ws.Range("D11").Value = -3.561242783948757
$ws.Range("E11").Value = 16.77205723344833
$ws.Range("F11").Value = 0.9985427644049873
$ws.Range("G11").Value = -0.7704248316864066
$ws.Range("H11").Value = -1.022015410261035
$ws.Range("I11").Value = 0.1356312661511834
$ws.Range("J11").Value = 4.353696125234475
$ws.Range("K11").Value = 213
$ws.Range("L11").Value = -30.49242943388767
$ws.Range("N11").Value = 4.356117816428284
$ws.Range("O11").Value = 5.022784483094951

